# Wetterdaten den Kraftwerken zuordnen
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New column P: "Trendquelle" (weather trend source) --------------------
# Header cell P1, formatted like the other header cells (A1, bold+yellow fill)
$ws.Range("A1").Copy() | Out-Null
$ws.Range("P1").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$ws.Range("P1").Value = "Trendquelle"

# Column P width
$ws.Columns.Item(16).ColumnWidth = 50.17

# --- Flag several existing power plants as "oNP allowed" (column J) --------
$ws.Range("J2").Value = 1
$ws.Range("J3").Value = 1
$ws.Range("J4").Value = 1
$ws.Range("J5").Value = 1
$ws.Range("J6").Value = 1
$ws.Range("J9").Value = 0
$ws.Range("J10").Value = 0

# --- Extend the shared formulas in H and I down through row 13 -------------
$ws.Range("H10:H13").Formula = "=4*C10"
$ws.Range("I10:I13").Formula = "=5*F10"

# --- Row 11: new wind power plant ------------------------------------------
$ws.Range("A11").Value = 10
$ws.Range("B11").Value = 1
$ws.Range("C11").Value = 100
$ws.Range("D11").Value = 0
$ws.Range("E11").Value = 1
$ws.Range("F11").Value = 0
$ws.Range("G11").Value = 3
$ws.Range("J11").Value = 0
$ws.Range("K11").Value = 0
$ws.Range("L11").Value = 0
$ws.Range("M11").Value = 0
$ws.Range("N11").Value = 0
$ws.Range("O11").Value = 0

# --- Row 12: new wind power plant -------------------------------------------
$ws.Range("A12").Value = 11
$ws.Range("B12").Value = 4
$ws.Range("C12").Value = 200
$ws.Range("D12").Value = 0
$ws.Range("E12").Value = 1
$ws.Range("F12").Value = 0
$ws.Range("G12").Value = 3
# N12/O12 did not have the grey-ish "s=3" formatting yet - match J12:M12
$ws.Range("J12").Copy() | Out-Null
$ws.Range("N12:O12").PasteSpecial(-4122) | Out-Null
$ws.Range("J12").Value = 0
$ws.Range("K12").Value = 0
$ws.Range("L12").Value = 0
$ws.Range("M12").Value = 0
$ws.Range("N12").Value = 0
$ws.Range("O12").Value = 0

# --- Row 13: new wind power plant -------------------------------------------
$ws.Range("A13").Value = 12
$ws.Range("B13").Value = 7
$ws.Range("C13").Value = 500
$ws.Range("D13").Value = 0
$ws.Range("E13").Value = 1
$ws.Range("F13").Value = 0
$ws.Range("G13").Value = 3
$ws.Range("J13").Value = 0
$ws.Range("K13").Value = 0
$ws.Range("L13").Value = 0
$ws.Range("M13").Value = 0
$ws.Range("N13").Value = 0
$ws.Range("O13").Value = 0

# --- Weather trend source file names, column P (order matters for the ------
#     shared-string table layout, matches how the workbook was authored)
$ws.Range("P13").Value = "weather/wind/Goteborg_Juli_2019.json"
$ws.Range("P11").Value = "weather/wind/Muenchen_Juli_2019.json"
$ws.Range("P12").Value = "weather/wind/Bremerhaven_Juli_2019.json"

# --- Remove the (now unused) formatted blank cell in row 16 ----------------
$ws.Range("P16").Clear()

# --- Update sheet view: selection/scroll position now on the new data -----
$ws.Range("P11").Select() | Out-Null
try {
    $excel.ActiveWindow.ScrollColumn = 14
    $excel.ActiveWindow.ScrollRow = 1
} catch {}
